# Updates Leve profit-tracking numbers ("currentAveragePrice" / NQ / HQ columns
# and their derived profit figures) across several crafting-job sheets, per the
# scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8160.7856
$ws.Range("J40").Value = 8031.375
$ws.Range("L40").Value = 8031.375
$ws.Range("N40").Value = -8381.375
$ws.Range("H74").Value = 6179.8
$ws.Range("I74").Value = 6179.8
$ws.Range("K74").Value = 6179.8
$ws.Range("M74").Value = -5243.8
$ws.Range("H77").Value = 6179.8
$ws.Range("I77").Value = 6179.8
$ws.Range("K77").Value = 30899
$ws.Range("M77").Value = -26219
$ws.Range("H100").Value = 3793
$ws.Range("I100").Value = 4500
$ws.Range("K100").Value = 4500
$ws.Range("M100").Value = -3959
$ws.Range("H101").Value = 671.6
$ws.Range("I101").Value = 671.6
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2014.8
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -392.8000000000002
$ws.Range("N101").ClearContents()
$ws.Range("H133").Value = 78526.61
$ws.Range("J133").Value = 78526.61
$ws.Range("L133").Value = 78526.61
$ws.Range("N133").Value = -88646.61
$ws.Range("H134").Value = 94216.664
$ws.Range("J134").Value = 94216.664
$ws.Range("L134").Value = 94216.664
$ws.Range("N134").Value = -104356.664
$ws.Range("H136").Value = 99995
$ws.Range("J136").Value = 99995
$ws.Range("L136").Value = 99995
$ws.Range("N136").Value = -110195
$ws.Range("H139").Value = 98406
$ws.Range("J139").Value = 98406
$ws.Range("L139").Value = 98406
$ws.Range("N139").Value = -108686
$ws.Range("H140").Value = 91989.42999999999
$ws.Range("J140").Value = 91989.42999999999
$ws.Range("L140").Value = 91989.42999999999
$ws.Range("N140").Value = -102349.43

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7380.382
$ws.Range("I32").Value = 3997.25
$ws.Range("K32").Value = 3997.25
$ws.Range("M32").Value = -3710.25
$ws.Range("H63").Value = 2298.1667
$ws.Range("I63").Value = 2317.8
$ws.Range("J63").Value = 2200
$ws.Range("K63").Value = 2317.8
$ws.Range("L63").Value = 2200
$ws.Range("M63").Value = -1631.8
$ws.Range("N63").Value = -3572
$ws.Range("H66").Value = 2298.1667
$ws.Range("I66").Value = 2317.8
$ws.Range("J66").Value = 2200
$ws.Range("K66").Value = 11589
$ws.Range("L66").Value = 11000
$ws.Range("M66").Value = -8157
$ws.Range("N66").Value = -17864
$ws.Range("H132").Value = 2351.6667
$ws.Range("I132").Value = 1963.3334
$ws.Range("K132").Value = 5890.0002
$ws.Range("M132").Value = -3360.0002

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1250258.8
$ws.Range("I4").Value = 1666963
$ws.Range("J4").Value = 146
$ws.Range("K4").Value = 1666963
$ws.Range("L4").Value = 146
$ws.Range("M4").Value = -1666848
$ws.Range("N4").Value = -376
$ws.Range("H22").Value = 14654304
$ws.Range("I22").Value = 14654304
$ws.Range("K22").Value = 14654304
$ws.Range("M22").Value = -14654131
$ws.Range("H82").Value = 18001.5
$ws.Range("I82").Value = 7601.8
$ws.Range("J82").Value = 70000
$ws.Range("K82").Value = 7601.8
$ws.Range("L82").Value = 70000
$ws.Range("M82").Value = -7218.8
$ws.Range("N82").Value = -70766
$ws.Range("H85").Value = 18001.5
$ws.Range("I85").Value = 7601.8
$ws.Range("J85").Value = 70000
$ws.Range("K85").Value = 7601.8
$ws.Range("L85").Value = 70000
$ws.Range("M85").Value = -6275.8
$ws.Range("N85").Value = -72652
$ws.Range("H94").Value = 1645.579
$ws.Range("I94").Value = 1363.8
$ws.Range("J94").Value = 2702.25
$ws.Range("K94").Value = 1363.8
$ws.Range("L94").Value = 2702.25
$ws.Range("M94").Value = -912.8
$ws.Range("N94").Value = -3604.25
$ws.Range("H132").Value = 45346.246
$ws.Range("J132").Value = 45346.246
$ws.Range("L132").Value = 45346.246
$ws.Range("N132").Value = -55466.246
$ws.Range("H138").Value = 89996.664
$ws.Range("J138").Value = 89996.664
$ws.Range("L138").Value = 89996.664
$ws.Range("N138").Value = -100276.664
$ws.Range("H140").Value = 72310.5
$ws.Range("J140").Value = 89990
$ws.Range("L140").Value = 89990
$ws.Range("N140").Value = -100350

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1237.25
$ws.Range("I22").Value = 1237.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1237.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -887.25
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 1644.0834
$ws.Range("I31").Value = 1406.2927
$ws.Range("J31").Value = 2157.2104
$ws.Range("K31").Value = 1406.2927
$ws.Range("L31").Value = 2157.2104
$ws.Range("M31").Value = -1111.2927
$ws.Range("N31").Value = -2747.2104
$ws.Range("H34").Value = 1644.0834
$ws.Range("I34").Value = 1406.2927
$ws.Range("J34").Value = 2157.2104
$ws.Range("K34").Value = 1406.2927
$ws.Range("L34").Value = 2157.2104
$ws.Range("M34").Value = -1204.2927
$ws.Range("N34").Value = -2561.2104
$ws.Range("H138").Value = 85839.836
$ws.Range("J138").Value = 89866
$ws.Range("L138").Value = 89866
$ws.Range("N138").Value = -100146

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 129.8
$ws.Range("J97").Value = 116.333336
$ws.Range("L97").Value = 349.000008
$ws.Range("N97").Value = -1341.000008

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20596.9
$ws.Range("J93").Value = 20596.9
$ws.Range("L93").Value = 20596.9
$ws.Range("N93").Value = -24340.9
$ws.Range("H132").Value = 3825.2083
$ws.Range("I132").Value = 2461.389
$ws.Range("J132").Value = 7916.6665
$ws.Range("K132").Value = 7384.167
$ws.Range("L132").Value = 23749.9995
$ws.Range("M132").Value = -4854.167
$ws.Range("N132").Value = -28809.9995
$ws.Range("H135").Value = 59146.1
$ws.Range("J135").Value = 59146.1
$ws.Range("L135").Value = 59146.1
$ws.Range("N135").Value = -69286.10000000001
$ws.Range("H140").Value = 98567.28999999999
$ws.Range("J140").Value = 98567.28999999999
$ws.Range("L140").Value = 98567.28999999999
$ws.Range("N140").Value = -108927.29

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4500.091
$ws.Range("I7").Value = 3093.1428
$ws.Range("K7").Value = 3093.1428
$ws.Range("M7").Value = -2981.1428
$ws.Range("H22").Value = 7377.75
$ws.Range("I22").Value = 1139.4
$ws.Range("J22").Value = 10213.363
$ws.Range("K22").Value = 1139.4
$ws.Range("L22").Value = 10213.363
$ws.Range("M22").Value = -844.4000000000001
$ws.Range("N22").Value = -10803.363
$ws.Range("H27").Value = 7377.75
$ws.Range("I27").Value = 1139.4
$ws.Range("J27").Value = 10213.363
$ws.Range("K27").Value = 1139.4
$ws.Range("L27").Value = 10213.363
$ws.Range("M27").Value = -1032.4
$ws.Range("N27").Value = -10427.363
$ws.Range("H40").Value = 2927816.8
$ws.Range("I40").Value = 4023.3845
$ws.Range("K40").Value = 4023.3845
$ws.Range("M40").Value = -3887.3845
$ws.Range("H46").Value = 7692.913
$ws.Range("I46").Value = 21669.6
$ws.Range("J46").Value = 3810.5
$ws.Range("K46").Value = 21669.6
$ws.Range("L46").Value = 3810.5
$ws.Range("M46").Value = -21481.6
$ws.Range("N46").Value = -4186.5
$ws.Range("H61").Value = 1657.421
$ws.Range("I61").Value = 1299.5
$ws.Range("J61").Value = 3566.3333
$ws.Range("K61").Value = 1299.5
$ws.Range("L61").Value = 3566.3333
$ws.Range("M61").Value = -1097.5
$ws.Range("N61").Value = -3970.3333
$ws.Range("H113").Value = 1657.421
$ws.Range("I113").Value = 1299.5
$ws.Range("J113").Value = 3566.3333
$ws.Range("K113").Value = 1299.5
$ws.Range("L113").Value = 3566.3333
$ws.Range("M113").Value = 870.5
$ws.Range("N113").Value = -7906.3333
$ws.Range("H126").Value = 4500.091
$ws.Range("I126").Value = 3093.1428
$ws.Range("K126").Value = 9279.428400000001
$ws.Range("M126").Value = -6809.428400000001
